$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H3").Value = "ZEW_Name"
